$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows below row 2 (was A1:A5, becomes A1:A2)
$ws.Range("A3:A5").ClearContents()

# Update the remaining data cell
$ws.Range("A2").Value = "TestProject001"

# Keep header text as-is (already "Projects") but make sure it matches
$ws.Range("A1").Value = "Projects"

# Update the active selection to A2
$ws.Range("A2").Select()
